$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are always plain text (coin name / URL) - safe to assign directly.
# Column E values keep their leading/trailing spaces so Excel never infers them as numbers.
# Column D sometimes looks like a genuine number ("211.60", "7.15", ...); assigning such a
# string via .Value would make Excel auto-convert the cell to a Number (like typing into the
# grid). The source data models these as text, so we force the cell to Text format first,
# write the literal string, then drop the temporary format back to the default "Normal"
# style so no stray formatting is left behind on the cell.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.739.18'
$ws.Range('E2').Value = '  +0.09%  '

Set-TextValue 'D3' '1.599.66'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  +0.24%  '

Set-TextValue 'D5' '211.60'
$ws.Range('E5').Value = '  +0.04%  '

Set-TextValue 'D6' '0.513'
$ws.Range('E6').Value = '  -0.03%  '

Set-TextValue 'D7' '1.01'
$ws.Range('E7').Value = '  +0.28%  '

Set-TextValue 'D8' '0.0619'
$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('E9').Value = '  -0.15%  '

Set-TextValue 'D10' '19.73'
$ws.Range('E10').Value = '  +1.11%  '

$ws.Range('E11').Value = '  +0.89%  '

Set-TextValue 'D12' '1.824.68'
$ws.Range('E12').Value = '  +0.04%  '

Set-TextValue 'D13' '1.619.01'
$ws.Range('E13').Value = '  +0.43%  '

Set-TextValue 'D14' '4.06'
$ws.Range('E14').Value = '  +0.56%  '

$ws.Range('E15').Value = '  +0.47%  '

Set-TextValue 'D16' '65.14'
$ws.Range('E16').Value = '  -0.29%  '

Set-TextValue 'D17' '26.722.46'
$ws.Range('E17').Value = '  +0.13%  '

Set-TextValue 'D18' '0.0₃0744'
$ws.Range('E18').Value = '  -2.00%  '

Set-TextValue 'D19' '209.81'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('E20').Value = '  +0.20%  '

Set-TextValue 'D21' '7.13'
$ws.Range('E21').Value = '  -0.25%  '

$ws.Range('E22').Value = '  +0.44%  '

$ws.Range('E23').Value = '  -2.09%  '

Set-TextValue 'D24' '9.02'
$ws.Range('E24').Value = '  +0.97%  '

Set-TextValue 'D25' '144.17'
$ws.Range('E25').Value = '  +0.74%  '

$ws.Range('E26').Value = '  +0.18%  '

Set-TextValue 'D27' '7.15'

$ws.Range('E28').Value = '  -0.51%  '

Set-TextValue 'D29' '15.35'
$ws.Range('E29').Value = '  +0.29%  '

$ws.Range('E30').Value = '  -1.92%  '

$ws.Range('E31').Value = '  +0.10%  '

$ws.Range('E32').Value = '  +0.66%  '

$ws.Range('E34').Value = '  +17.57%  '

Set-TextValue 'D35' '1.277.73'
$ws.Range('E35').Value = '  -0.95%  '

$ws.Range('E36').Value = '  +0.98%  '

$ws.Range('E37').Value = '  -0.41%  '

$ws.Range('E38').Value = '  -3.58%  '

$ws.Range('E39').Value = '  -1.36%  '

$ws.Range('E40').Value = '  +0.03%  '

$ws.Range('E41').Value = '  +2.33%  '

Set-TextValue 'D42' '5.45'
$ws.Range('E42').Value = '  +0.28%  '

Set-TextValue 'D44' '62.67'
$ws.Range('E44').Value = '  -0.63%  '

Set-TextValue 'D45' '1.735.90'
$ws.Range('E45').Value = '  -0.07%  '

Set-TextValue 'D46' '90.51'
$ws.Range('E46').Value = '  -0.93%  '

$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('E48').Value = '  +2.63%  '

$ws.Range('E49').Value = '  +0.82%  '

Set-TextValue 'D50' '7.58'
$ws.Range('E50').Value = '  +2.79%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D51' '1.01'
$ws.Range('E51').Value = '  +0.28%  '
